$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from SCD0195 to SCD0011
$ws.Name = "SCD0011"

# Update cell B2 value from "DGS-210" to "SCD0011-026"
$ws.Range("B2").Value = "SCD0011-026"

# Adjust column B width (closest achievable value given Excel's pixel-quantized
# column-width storage; COM ColumnWidth values in this neighborhood all collapse
# to the same stored width, which is the nearest representable value to 8.85546875)
$ws.Range("B:B").ColumnWidth = 8.02

# Update selection to B3 (also resets the scrolled top-left cell back to default)
$ws.Range("B3").Select()
